# fix(cad): update ecad / mcad
# The "Alternate Part Number (PCBWay Preferred)" values in column H were
# entered one row too low. Shift them up by one row so each alternate
# part number lines up with the correct BOM line:
#   H19 (MCP73831/OT, IC2)            -> MCP73831T-2ACI/OT
#   H20 (MICRO-SD_CARD_SOCKET, J2)    -> MEM2051-00-195-00-A
#   H21 (PINHD-2X2REVERSED, JP3)      -> (no alternate part number)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("H19").Value = "MCP73831T-2ACI/OT"
$ws.Range("H20").Value = "MEM2051-00-195-00-A"
$ws.Range("H21").Value = ""

# Update the sheet's last-saved cell selection to reflect where the user
# was working when they made this fix.
$ws.Range("H20").Select()
